$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
# A1 stays "Ano". B1/C1 change meaning, D1:L1 are new columns.
$ws.Range("B1").Value = "Produtos"
$ws.Range("C1").Value = "Serviços"

$headers = @(
    "Total",
    "Serviços (%)",
    "Produtos (%)",
    "Evolução Serviços (%)",
    "Evolução Produtos (%)",
    "Evolução Total (%)",
    "Qtd Produtos",
    "Qtd Serviços",
    "Total Itens"
)
$col = 4
foreach ($h in $headers) {
    $ws.Cells.Item(1, $col).Value = $h
    $col = $col + 1
}

# Copy the header style (bold, centered, bordered) from A1 onto the new D1:L1 headers
$ws.Range("A1").Copy()
$ws.Range("D1:L1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Data rows (2-9) ---
$ws.Cells.Item(2, 1).Value = 2018
$ws.Cells.Item(2, 2).Value = 484712.74
$ws.Cells.Item(2, 3).Value = 234945
$ws.Cells.Item(2, 4).Value = 719657.74
$ws.Cells.Item(2, 5).Value = 32.64676900438811
$ws.Cells.Item(2, 6).Value = 67.3532309956119
$ws.Cells.Item(2, 10).Value = 1780
$ws.Cells.Item(2, 11).Value = 1351
$ws.Cells.Item(2, 12).Value = 3131

$ws.Cells.Item(3, 1).Value = 2019
$ws.Cells.Item(3, 2).Value = 1161517.53
$ws.Cells.Item(3, 3).Value = 635735.5600000001
$ws.Cells.Item(3, 4).Value = 1797253.09
$ws.Cells.Item(3, 5).Value = 35.37262300659057
$ws.Cells.Item(3, 6).Value = 64.62737699340943
$ws.Cells.Item(3, 7).Value = 170.5890995764967
$ws.Cells.Item(3, 8).Value = 139.6300806948049
$ws.Cells.Item(3, 9).Value = 149.7372000751357
$ws.Cells.Item(3, 10).Value = 19043
$ws.Cells.Item(3, 11).Value = 3095
$ws.Cells.Item(3, 12).Value = 22138

$ws.Cells.Item(4, 1).Value = 2020
$ws.Cells.Item(4, 2).Value = 2259471.49
$ws.Cells.Item(4, 3).Value = 1114043
$ws.Cells.Item(4, 4).Value = 3373514.49
$ws.Cells.Item(4, 5).Value = 33.02321668699872
$ws.Cells.Item(4, 6).Value = 66.97678331300128
$ws.Cells.Item(4, 7).Value = 75.23685477024441
$ws.Cells.Item(4, 8).Value = 94.52754105226464
$ws.Cells.Item(4, 9).Value = 87.70391931834153
$ws.Cells.Item(4, 10).Value = 42777
$ws.Cells.Item(4, 11).Value = 4019
$ws.Cells.Item(4, 12).Value = 46796

$ws.Cells.Item(5, 1).Value = 2021
$ws.Cells.Item(5, 2).Value = 5399571.25
$ws.Cells.Item(5, 3).Value = 1275732.01
$ws.Cells.Item(5, 4).Value = 6675303.26
$ws.Cells.Item(5, 5).Value = 19.11122177241727
$ws.Cells.Item(5, 6).Value = 80.88877822758272
$ws.Cells.Item(5, 7).Value = 14.51371356401863
$ws.Cells.Item(5, 8).Value = 138.9749671061351
$ws.Cells.Item(5, 9).Value = 97.87385765756707
$ws.Cells.Item(5, 10).Value = 64210
$ws.Cells.Item(5, 11).Value = 3896
$ws.Cells.Item(5, 12).Value = 68106

$ws.Cells.Item(6, 1).Value = 2022
$ws.Cells.Item(6, 2).Value = 6143248.8
$ws.Cells.Item(6, 3).Value = 1358552.5
$ws.Cells.Item(6, 4).Value = 7501801.3
$ws.Cells.Item(6, 5).Value = 18.10968387019262
$ws.Cells.Item(6, 6).Value = 81.89031612980739
$ws.Cells.Item(6, 7).Value = 6.49199748464413
$ws.Cells.Item(6, 8).Value = 13.77290002423803
$ws.Cells.Item(6, 9).Value = 12.38143059286267
$ws.Cells.Item(6, 10).Value = 67109
$ws.Cells.Item(6, 11).Value = 3959
$ws.Cells.Item(6, 12).Value = 71068

$ws.Cells.Item(7, 1).Value = 2023
$ws.Cells.Item(7, 2).Value = 5961894
$ws.Cells.Item(7, 3).Value = 1198106.18
$ws.Cells.Item(7, 4).Value = 7160000.18
$ws.Cells.Item(7, 5).Value = 16.73332611564264
$ws.Cells.Item(7, 6).Value = 83.26667388435737
$ws.Cells.Item(7, 7).Value = -11.81009346344731
$ws.Cells.Item(7, 8).Value = -2.952099221506377
$ws.Cells.Item(7, 9).Value = -4.556253975961743
$ws.Cells.Item(7, 10).Value = 70506
$ws.Cells.Item(7, 11).Value = 3843
$ws.Cells.Item(7, 12).Value = 74349

$ws.Cells.Item(8, 1).Value = 2024
$ws.Cells.Item(8, 2).Value = 8510815.33
$ws.Cells.Item(8, 3).Value = 1052835.53
$ws.Cells.Item(8, 4).Value = 9563650.859999999
$ws.Cells.Item(8, 5).Value = 11.00871984362675
$ws.Cells.Item(8, 6).Value = 88.99128015637326
$ws.Cells.Item(8, 7).Value = -12.12502300922944
$ws.Cells.Item(8, 8).Value = 42.75354996247837
$ws.Cells.Item(8, 9).Value = 33.57053937951158
$ws.Cells.Item(8, 10).Value = 73484
$ws.Cells.Item(8, 11).Value = 3279
$ws.Cells.Item(8, 12).Value = 76763

$ws.Cells.Item(9, 1).Value = 2025
$ws.Cells.Item(9, 2).Value = 3098380.93
$ws.Cells.Item(9, 3).Value = 322046.99
$ws.Cells.Item(9, 4).Value = 3420427.92
$ws.Cells.Item(9, 5).Value = 9.415400573621795
$ws.Cells.Item(9, 6).Value = 90.5845994263782
$ws.Cells.Item(9, 7).Value = -69.4114625861838
$ws.Cells.Item(9, 8).Value = -63.59478134746463
$ws.Cells.Item(9, 9).Value = -64.23512348923201
$ws.Cells.Item(9, 10).Value = 22866
$ws.Cells.Item(9, 11).Value = 908
$ws.Cells.Item(9, 12).Value = 23774
